# Trade #103 closed at 2026-02-18 00:37:32 - unknown UNKNOWN +0.000%
#
# This closes the open MarketMaking trade (row 132 on "All Trades" / row 52
# on "MarketMaking") with an early_exit, and records a brand-new open
# HighProbConvergence trade (row 161 on "All Trades" / row 23 on
# "HighProbConvergence"). The Summary and Strategy Status roll-up sheets are
# updated to reflect the newly closed trade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - portfolio level roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.07   # Current Capital
$summary.Range("B4").Value = 0.18      # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 131       # Total Trades
$summary.Range("B7").Value = 62        # Winning Trades
$summary.Range("B9").Value = 47.33     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking strategy row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.15000000000001
$status.Range("D6").Value = 51
$status.Range("E6").Value = -0.66
$status.Range("F6").Value = -0.85
$status.Range("G6").Value = 47.06

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close the existing MarketMaking trade recorded in row 132
$allTrades.Range("G132").Value = 0.96
$allTrades.Range("H132").Value = "CLOSED"
$allTrades.Range("I132").Value = 2.1277
$allTrades.Range("J132").Value = 0.02
$allTrades.Range("K132").Value = 99.15000000000001
$allTrades.Range("L132").Value = "early_exit"
$allTrades.Range("M132").Value = 0.14

# Append the newly opened HighProbConvergence trade as row 161
$allTrades.Range("A161").Value = 160
# "2026-02-18" looks like a date to Excel's input parser, so it would get
# silently converted to a date serial number; force it to stay literal text
# (matching every other date cell in this workbook, which are plain strings)
# by temporarily marking the cell as Text, then restore the default format.
$allTrades.Range("B161").NumberFormat = "@"
$allTrades.Range("B161").Value = "2026-02-18"
$allTrades.Range("B161").ClearFormats()
$allTrades.Range("C161").Value = "00:37:26"
$allTrades.Range("D161").Value = "HighProbConvergence"
$allTrades.Range("E161").Value = "UP"
$allTrades.Range("F161").Value = 0.9399999999999999
$allTrades.Range("H161").Value = "OPEN"
$allTrades.Range("I161").Value = 0
$allTrades.Range("J161").Value = 0
$allTrades.Range("K161").Value = 100.4130057263667
$allTrades.Range("M161").Value = 0
$allTrades.Range("N161").Value = 0
$allTrades.Range("O161").Value = 0
$allTrades.Range("P161").Value = 0.95
$allTrades.Range("Q161").Value = "Mean reversion UP: price 1.85% below mean (z=-4.36)"

# ---------------------------------------------------------------------
# HighProbConvergence sheet - append the same new trade as row 23
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("A23").Value = 160
$hpc.Range("B23").NumberFormat = "@"
$hpc.Range("B23").Value = "2026-02-18"
$hpc.Range("B23").ClearFormats()
$hpc.Range("C23").Value = "00:37:26"
$hpc.Range("D23").Value = "HighProbConvergence"
$hpc.Range("E23").Value = "UP"
$hpc.Range("F23").Value = 0.9399999999999999
$hpc.Range("H23").Value = "OPEN"
$hpc.Range("I23").Value = 0
$hpc.Range("J23").Value = 0
$hpc.Range("K23").Value = 100.4130057263667
$hpc.Range("L23").Value = 0
$hpc.Range("M23").Value = 0
$hpc.Range("N23").Value = 0.95
$hpc.Range("O23").Value = "Mean reversion UP: price 1.85% below mean (z=-4.36)"
$hpc.Range("Q23").Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - close the same trade recorded in row 52
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G52").Value = 0.96
$mm.Range("H52").Value = "CLOSED"
$mm.Range("I52").Value = 2.1277
$mm.Range("J52").Value = 0.02
$mm.Range("K52").Value = 99.15000000000001
$mm.Range("P52").Value = "early_exit"
$mm.Range("Q52").Value = 0.14
